$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.95
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 2.18
$ws.Range("J2").Value = 4.4
$ws.Range("K2").Value = 1.9
$ws.Range("L2").Value = 2.8
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 5.3
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.4
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.47
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.42
$ws.Range("U2").Value = 2.02
$ws.Range("V2").Value = 1.7
$ws.Range("AB2").Value = 55
$ws.Range("AC2").Value = 5.3
$ws.Range("AD2").Value = 5.5
$ws.Range("AE2").Value = 16
$ws.Range("AF2").Value = 100
$ws.Range("AG2").Value = 900
$ws.Range("AH2").Value = 5.6
$ws.Range("AI2").Value = 9.25
$ws.Range("AJ2").Value = 9
$ws.Range("AK2").Value = 22
$ws.Range("AL2").Value = 22
$ws.Range("AM2").Value = 40
$ws.Range("AN2").Value = 5.6
$ws.Range("AQ2").Value = 150
$ws.Range("AT2").Value = 2.37
$ws.Range("AU2").Value = 7.1
$ws.Range("AV2").Value = 70
$ws.Range("AW2").Value = 3.9
$ws.Range("AX2").Value = 11.75
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 50
$ws.Range("BA2").Value = 90
$ws.Range("BB2").Value = 300
$ws.Range("H3").Value = 2.7
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 2.95
$ws.Range("K3").Value = 1.85
$ws.Range("L3").Value = 4.4
$ws.Range("N3").Value = 5
$ws.Range("P3").Value = 2.32
$ws.Range("Y3").Value = 9.5
$ws.Range("AA3").Value = 24
$ws.Range("AB3").Value = 45
$ws.Range("AC3").Value = 5
$ws.Range("AD3").Value = 5.4
$ws.Range("AE3").Value = 17.5
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 60
$ws.Range("AN3").Value = 3.95
$ws.Range("AP3").Value = 24
$ws.Range("AQ3").Value = 60
$ws.Range("AR3").Value = 110
$ws.Range("AS3").Value = 400
$ws.Range("AU3").Value = 7.5
$ws.Range("AW3").Value = 5.4
$ws.Range("AX3").Value = 23
$ws.Range("AY3").Value = 32
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 1.83
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 7
$ws.Range("Z7").Value = 9
$ws.Range("AD7").Value = 10
$ws.Range("AF7").Value = 67
$ws.Range("AK7").Value = 81
$ws.Range("AO7").Value = 6.5
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 17
$ws.Range("AZ7").Value = 151
